# Automatische test-sync: 2025-06-29 15:04:50
# Append the 11th test-mail log entry (row 26) to the "Logs" sheet and
# bump the "Bestelling / Levering" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

$row = 26

$logs.Range("A$row").Value = "Ik heb geen orderbevestiging gekregen"
$logs.Range("B$row").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$row").Value = "Testmail #11: Ik heb geen orderbevestiging gekregen"
$logs.Range("D$row").Value = "Bestelling / Levering"
$logs.Range("E$row").Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om u beter van dienst te kunnen zijn, zou u ons uw bestelnummer of de datum van uw bestelling kunnen doorgeven? Hiermee kunnen we nagaan wat er mogelijk is misgegaan met de orderbevestiging.`nWe doen ons best om dit zo snel mogelijk voor u op te lossen. Alvast bedankt voor uw medewerking!`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F$row").Value = "2025-06-29 15:04:32"
$logs.Range("G$row").Value = "Ja"
$logs.Range("H$row").Value = "Nee"
$logs.Range("I$row").Value = "Ja"

# The multi-line Antwoord text auto-marks the row with a custom height;
# AutoFit() puts it back to the sheet's (non-custom) default so row 26
# looks like every other data row.
$logs.Rows.Item($row).AutoFit()

# Extend the four conditional-formatting blocks (D, G, H, I columns) so
# they keep covering the data through the newly added row.
$fcD = $logs.Range("D2:D25").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($logs.Range("D2:D26"))

$fcG = $logs.Range("G2:G25").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($logs.Range("G2:G26"))

$fcH = $logs.Range("H2:H25").FormatConditions
$fcH.Item(1).ModifyAppliesToRange($logs.Range("H2:H26"))

$fcI = $logs.Range("I2:I25").FormatConditions
$fcI.Item(1).ModifyAppliesToRange($logs.Range("I2:I26"))

# Update the Dashboard summary count for "Bestelling / Levering" (6 -> 7).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 7
